$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) storage for numeric-looking price cells so Excel
# does not auto-convert them to floating point numbers, matching the
# original inline-string representation.
$textCells = @(
    'D5',
    'D10',
    'D11',
    'D12',
    'D15',
    'D17',
    'D21',
    'D22',
    'D25',
    'D31',
    'D32',
    'D36',
    'D42',
    'D44',
    'D46',
    'D49',
    'D50',
    'D51'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values (coin metadata, prices, and hourly volume %).
$ws.Range('D2').Value = '26.282.17'
$ws.Range('E2').Value = '  +1.93%  '
$ws.Range('D3').Value = '1.647.82'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('D5').Value = '216.98'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('D10').Value = '19.94'
$ws.Range('E10').Value = '  +1.57%  '
$ws.Range('D11').Value = '0.0794'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '4.30'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').Value = '1.874.56'
$ws.Range('E13').Value = '  +0.54%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.665.22'
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('D15').Value = '0.548'
$ws.Range('E15').Value = '  -2.67%  '
$ws.Range('D16').Value = '0.0₃0767'
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').Value = '63.37'
$ws.Range('D18').Value = '26.273.59'
$ws.Range('E18').Value = '  +1.71%  '
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('E20').Value = '  -0.61%  '
$ws.Range('D21').Value = '195.69'
$ws.Range('E21').Value = '  +1.67%  '
$ws.Range('D22').Value = '10.09'
$ws.Range('E22').Value = '  +1.29%  '
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('E24').Value = '  -3.41%  '
$ws.Range('D25').Value = '143.49'
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E27').Value = '  +0.88%  '
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('E29').Value = '  +0.99%  '
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').Value = '0.0505'
$ws.Range('E31').Value = '  +2.45%  '
$ws.Range('D32').Value = '3.36'
$ws.Range('E32').Value = '  +0.94%  '
$ws.Range('E33').Value = '  +0.92%  '
$ws.Range('E34').Value = '  +2.15%  '
$ws.Range('E35').Value = '  +1.27%  '
$ws.Range('D36').Value = '0.914'
$ws.Range('E36').Value = '  +1.00%  '
$ws.Range('D37').Value = '1.139.33'
$ws.Range('E37').Value = '  +0.28%  '
$ws.Range('E38').Value = '  +1.57%  '
$ws.Range('E39').Value = '  -1.75%  '
$ws.Range('E40').Value = '  +1.26%  '
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('D42').Value = '100.56'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('D44').Value = '0.800'
$ws.Range('E44').Value = '  -0.51%  '
$ws.Range('D45').Value = '1.783.27'
$ws.Range('E45').Value = '  +0.55%  '
$ws.Range('D46').Value = '57.16'
$ws.Range('E46').Value = '  +3.36%  '
$ws.Range('E47').Value = '  +4.26%  '
$ws.Range('E48').Value = '  +3.04%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '7.72'
$ws.Range('E49').Value = '  +3.53%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '0.418'
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('D51').Value = '0.0971'
$ws.Range('E51').Value = '  +1.72%  '
